# Auto-generated edit script: update ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# with refreshed market-price + profit figures (scheduled runner update).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3621.889
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3621.889
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3621.889
$ws.Range("N40").Value = -3971.889
$ws.Range("M40").ClearContents()

$ws.Range("H64").Value = 8440
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 8440
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 8440
$ws.Range("N64").Value = -8936
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 8440
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 8440
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 8440
$ws.Range("N67").Value = -10156
$ws.Range("M67").ClearContents()

$ws.Range("H70").Value = 1691381.4
$ws.Range("I70").Value = 5062772
$ws.Range("J70").Value = 5686
$ws.Range("K70").Value = 15188316
$ws.Range("L70").Value = 17058
$ws.Range("M70").Value = -15188046
$ws.Range("N70").Value = -17598

$ws.Range("H73").Value = 1691381.4
$ws.Range("I73").Value = 5062772
$ws.Range("J73").Value = 5686
$ws.Range("K73").Value = 15188316
$ws.Range("L73").Value = 17058
$ws.Range("M73").Value = -15187380
$ws.Range("N73").Value = -18930

$ws.Range("H97").Value = 1200
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 3600
$ws.Range("N97").Value = -4592

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3602.6
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 3602.6
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 3602.6
$ws.Range("N2").Value = -3828.6
$ws.Range("M2").ClearContents()

$ws.Range("H45").Value = 3449.25
$ws.Range("I45").Value = 3449.25
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 3449.25
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -3072.25

$ws.Range("H88").Value = 6399.75
$ws.Range("I88").Value = 1100
$ws.Range("J88").Value = 7156.857
$ws.Range("K88").Value = 1100
$ws.Range("L88").Value = 7156.857
$ws.Range("M88").Value = -694
$ws.Range("N88").Value = -7968.857

$ws.Range("H91").Value = 6399.75
$ws.Range("I91").Value = 1100
$ws.Range("J91").Value = 7156.857
$ws.Range("K91").Value = 1100
$ws.Range("L91").Value = 7156.857
$ws.Range("M91").Value = 304
$ws.Range("N91").Value = -9964.857

$ws.Range("H103").Value = 80000
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 80000
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 80000
$ws.Range("N103").Value = -82344

$ws.Range("H116").Value = 3602.6
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 3602.6
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 3602.6
$ws.Range("N116").Value = -8190.6
$ws.Range("M116").ClearContents()

$ws.Range("H132").Value = 5779
$ws.Range("I132").Value = 5779
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 17337
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -14807

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 534.125
$ws.Range("I10").Value = 259.6
$ws.Range("J10").Value = 991.6667
$ws.Range("K10").Value = 259.6
$ws.Range("L10").Value = 991.6667
$ws.Range("M10").Value = -120.6
$ws.Range("N10").Value = -1269.6667

$ws.Range("H22").Value = 925
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -650

$ws.Range("H31").Value = 3052.7273
$ws.Range("I31").Value = 1614.4286
$ws.Range("J31").Value = 5569.75
$ws.Range("K31").Value = 1614.4286
$ws.Range("L31").Value = 5569.75
$ws.Range("M31").Value = -1319.4286

$ws.Range("H34").Value = 3052.7273
$ws.Range("I34").Value = 1614.4286
$ws.Range("J34").Value = 5569.75
$ws.Range("K34").Value = 1614.4286
$ws.Range("L34").Value = 5569.75
$ws.Range("M34").Value = -1412.4286

$ws.Range("H86").Value = 12500
$ws.Range("I86").Value = 10000
$ws.Range("J86").Value = 15000
$ws.Range("K86").Value = 10000
$ws.Range("L86").Value = 15000
$ws.Range("M86").Value = -8877
$ws.Range("N86").Value = -17246

$ws.Range("H89").Value = 12500
$ws.Range("I89").Value = 10000
$ws.Range("J89").Value = 15000
$ws.Range("K89").Value = 50000
$ws.Range("L89").Value = 75000
$ws.Range("M89").Value = -44384
$ws.Range("N89").Value = -86232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 148.375
$ws.Range("I2").Value = 174.57143
$ws.Range("J2").Value = 111.7
$ws.Range("K2").Value = 1047.42858
$ws.Range("L2").Value = 670.2
$ws.Range("M2").Value = -934.42858
$ws.Range("N2").Value = -896.2

$ws.Range("H34").Value = 2330
$ws.Range("I34").Value = 353
$ws.Range("J34").Value = 4966
$ws.Range("K34").Value = 1059
$ws.Range("L34").Value = 14898
$ws.Range("M34").Value = -975
$ws.Range("N34").Value = -15066

$ws.Range("H39").Value = 3339.2
$ws.Range("I39").Value = 1349.5
$ws.Range("J39").Value = 4665.6665
$ws.Range("K39").Value = 4048.5
$ws.Range("L39").Value = 13996.9995
$ws.Range("M39").Value = -3754.5
$ws.Range("N39").Value = -14584.9995

$ws.Range("H55").Value = 224.5
$ws.Range("I55").Value = 224.5
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 673.5
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -496.5

$ws.Range("H107").Value = 1759.875
$ws.Range("I107").Value = 1276
$ws.Range("J107").Value = 1829
$ws.Range("K107").Value = 3828
$ws.Range("L107").Value = 5487
$ws.Range("M107").Value = -1908
$ws.Range("N107").Value = -9327

$ws.Range("H137").Value = 8592.6
$ws.Range("I137").Value = 2000
$ws.Range("J137").Value = 10240.75
$ws.Range("K137").Value = 6000
$ws.Range("L137").Value = 30722.25
$ws.Range("M137").Value = -900
$ws.Range("N137").Value = -40922.25

$ws.Range("H138").Value = 1099
$ws.Range("I138").Value = 648.5
$ws.Range("J138").Value = 2000
$ws.Range("K138").Value = 1945.5
$ws.Range("L138").Value = 6000
$ws.Range("M138").Value = 3194.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8156.9165
$ws.Range("I70").Value = 7276.6
$ws.Range("J70").Value = 8785.714
$ws.Range("K70").Value = 7276.6
$ws.Range("L70").Value = 8785.714
$ws.Range("M70").Value = -7006.6

$ws.Range("H73").Value = 8156.9165
$ws.Range("I73").Value = 7276.6
$ws.Range("J73").Value = 8785.714
$ws.Range("K73").Value = 7276.6
$ws.Range("L73").Value = 8785.714
$ws.Range("M73").Value = -6340.6

$ws.Range("H99").Value = 33670.43
$ws.Range("I99").Value = 45235.75
$ws.Range("J99").Value = 18250
$ws.Range("K99").Value = 45235.75
$ws.Range("L99").Value = 18250
$ws.Range("M99").Value = -42989.75
$ws.Range("N99").Value = -22742

$ws.Range("H111").Value = 60000
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 60000
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 60000
$ws.Range("N111").Value = -66134

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 653
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 653
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 653
$ws.Range("N20").Value = -1105

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H26").Value = 2394
$ws.Range("I26").Value = 2859
$ws.Range("J26").Value = 999
$ws.Range("K26").Value = 2859
$ws.Range("L26").Value = 999
$ws.Range("M26").Value = -2564
$ws.Range("N26").Value = -1589

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()

$ws.Range("H61").Value = 3058
$ws.Range("I61").Value = 3058
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3058
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2856
$ws.Range("N61").ClearContents()

$ws.Range("H105").Value = 62450
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 62450
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 62450
$ws.Range("N105").Value = -69438

$ws.Range("H110").Value = 29990
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 29990
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 29990
$ws.Range("N110").Value = -38170

$ws.Range("H113").Value = 3058
$ws.Range("I113").Value = 3058
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3058
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -888
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 1497.5
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1497.5
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 1497.5
$ws.Range("N8").Value = -1777.5
$ws.Range("M8").ClearContents()

$ws.Range("H69").Value = 5939.6665
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 5939.6665
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 5939.6665
$ws.Range("N69").Value = -7437.6665
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 5939.6665
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 5939.6665
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 17818.9995
$ws.Range("N72").Value = -25306.9995
$ws.Range("M72").ClearContents()

$ws.Range("H97").Value = 54996.5
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 54996.5
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 54996.5
$ws.Range("N97").Value = -56978.5

$ws.Range("H126").Value = 4632.3335
$ws.Range("I126").Value = 4632.3335
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 13897.0005
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -11427.0005
